$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-04-11 Thursday" "2024-04-12 Friday"

Replace-Text "972÷9=108, 0" "574÷9=63, 7"
Replace-Text "424÷3=141, 1" "927÷9=103, 0"
Replace-Text "659÷3=219, 2" "622÷6=103, 4"
Replace-Text "113÷8=14, 1" "502÷5=100, 2"
Replace-Text "589÷9=65, 4" "968÷6=161, 2"

Replace-Text "826÷3=275, 1" "673÷6=112, 1"
Replace-Text "806÷7=115, 1" "157÷6=26, 1"
Replace-Text "394÷4=98, 2" "676÷4=169, 0"
Replace-Text "801÷6=133, 3" "390÷4=97, 2"
Replace-Text "778÷3=259, 1" "486÷4=121, 2"

Replace-Text "765÷2=382, 1" "338÷2=169, 0"
Replace-Text "463÷7=66, 1" "784÷3=261, 1"
Replace-Text "356÷9=39, 5" "234÷8=29, 2"
Replace-Text "307÷7=43, 6" "740÷4=185, 0"
Replace-Text "638÷3=212, 2" "188÷6=31, 2"

Replace-Text "577÷7=82, 3" "740÷2=370, 0"
Replace-Text "603÷6=100, 3" "942÷6=157, 0"
Replace-Text "982÷2=491, 0" "362÷6=60, 2"
Replace-Text "102÷5=20, 2" "717÷8=89, 5"
Replace-Text "494÷3=164, 2" "376÷9=41, 7"

Replace-Text "420÷2=210, 0" "478÷8=59, 6"
Replace-Text "131÷5=26, 1" "236÷5=47, 1"
Replace-Text "566÷2=283, 0" "617÷6=102, 5"
Replace-Text "208÷8=26, 0" "796÷3=265, 1"
Replace-Text "247÷8=30, 7" "439÷6=73, 1"
